# select_school_surveys.xlsx -- additions to select school data file
#
# 1. Adds a new column AA "teachers_recommend_school_to_families" with
#    values for the existing 8 schools (rows 2-9).
# 2. Adds three new school survey rows (10-12: q151/woodside,
#    q148/east-elmhurst, q361/woodside) across all columns A:AA.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New column AA: copy the formatting from column Z (so the new cells
#    pick up the same style index instead of creating new style/font
#    entries), then write the header and the per-school values.
# ---------------------------------------------------------------------
$ws.Range("Z1:Z9").Copy()
$ws.Range("AA1:AA9").PasteSpecial(-4122)

$ws.Range("AA1").Value = "teachers_recommend_school_to_families"

# ---------------------------------------------------------------------
# 2. New rows 10-12: copy the formatting from row 9 (now that it spans
#    A:AA) down into the new rows first ...
# ---------------------------------------------------------------------
$ws.Range("A9:AA9").Copy()
$ws.Range("A10:AA12").PasteSpecial(-4122)

# ... then fill in the school / neighborhood text columns ...
$ws.Range("A10").Value = "q151"
$ws.Range("A11").Value = "q148"
$ws.Range("A12").Value = "q361"

$ws.Range("B10").Value = "woodside"
$ws.Range("B11").Value = "east-elmhurst"
$ws.Range("B12").Value = "woodside"

# ... and finally the numeric survey columns for every new row.
$ws.Range("AA2").Value = 82
$ws.Range("AA3").Value = 95
$ws.Range("AA4").Value = 100
$ws.Range("AA5").Value = 83
$ws.Range("AA6").Value = 89
$ws.Range("AA7").Value = 89
$ws.Range("AA8").Value = 98
$ws.Range("AA9").Value = 82

$ws.Range("C10").Value = 30
$ws.Range("D10").Value = 68
$ws.Range("E10").Value = 30
$ws.Range("F10").Value = 84
$ws.Range("G10").Value = 89
$ws.Range("H10").Value = 99
$ws.Range("I10").Value = 100
$ws.Range("J10").Value = 97
$ws.Range("K10").Value = 100
$ws.Range("L10").Value = 98
$ws.Range("M10").Value = 99
$ws.Range("N10").Value = 97
$ws.Range("O10").Value = 98
$ws.Range("P10").Value = 100
$ws.Range("Q10").Value = 99
$ws.Range("R10").Value = 97
$ws.Range("S10").Value = 95
$ws.Range("T10").Value = 97
$ws.Range("U10").Value = 98
$ws.Range("V10").Value = 100
$ws.Range("W10").Value = 93
$ws.Range("X10").Value = 89
$ws.Range("Y10").Value = 98
$ws.Range("Z10").Value = 100
$ws.Range("AA10").Value = 100

$ws.Range("C11").Value = 30
$ws.Range("D11").Value = 52
$ws.Range("E11").Value = 33
$ws.Range("F11").Value = 76
$ws.Range("G11").Value = 73
$ws.Range("H11").Value = 91
$ws.Range("I11").Value = 91
$ws.Range("J11").Value = 74
$ws.Range("K11").Value = 84
$ws.Range("L11").Value = 89
$ws.Range("M11").Value = 86
$ws.Range("N11").Value = 96
$ws.Range("O11").Value = 97
$ws.Range("P11").Value = 90
$ws.Range("Q11").Value = 87
$ws.Range("R11").Value = 92
$ws.Range("S11").Value = 92
$ws.Range("T11").Value = 97
$ws.Range("U11").Value = 90
$ws.Range("V11").Value = 75
$ws.Range("W11").Value = 81
$ws.Range("X11").Value = 87
$ws.Range("Y11").Value = 88
$ws.Range("Z11").Value = 98
$ws.Range("AA11").Value = 85

$ws.Range("C12").Value = 30
$ws.Range("D12").Value = 100
$ws.Range("E12").Value = 81
$ws.Range("F12").Value = 84
$ws.Range("G12").Value = 83
$ws.Range("H12").Value = 99
$ws.Range("I12").Value = 97
$ws.Range("J12").Value = 99
$ws.Range("K12").Value = 98
$ws.Range("L12").Value = 89
$ws.Range("M12").Value = 95
$ws.Range("N12").Value = 98
$ws.Range("O12").Value = 98
$ws.Range("P12").Value = 90
$ws.Range("Q12").Value = 99
$ws.Range("R12").Value = 98
$ws.Range("S12").Value = 96
$ws.Range("T12").Value = 96
$ws.Range("U12").Value = 94
$ws.Range("V12").Value = 94
$ws.Range("W12").Value = 85
$ws.Range("X12").Value = 87
$ws.Range("Y12").Value = 98
$ws.Range("Z12").Value = 100
$ws.Range("AA12").Value = 97

# ---------------------------------------------------------------------
# 3. View cosmetics that changed along with the data edit: zoom level
#    and the active selection.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 175
$ws.Range("U4").Select() | Out-Null
